$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Code Quality")
$ws.Activate()

$ws.Range("C4").Value = 80.8
$ws.Range("C5").Value = 71.6
$ws.Range("C6").Value = 1

$ws.Range("D12").Select()
